$wb = $excel.ActiveWorkbook

# --- Fix view/selection state on existing sheets to match target ---
$sheetCommon = $wb.Worksheets.Item("common")
$sheetCommon.Activate()
$sheetCommon.Range("H26").Select()

$sheetGroup = $wb.Worksheets.Item("group-page.html")
$sheetGroup.Activate()
$sheetGroup.Range("D1").Select()

# --- Add the new "creating-pairs-list-page.html" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "creating-pairs-list-page.html"

# --- Column widths (approximate OOXML stored widths of 25.83 / 24.51 / 34.87 / 40.83) ---
$newSheet.Columns.Item(1).ColumnWidth = 25
$newSheet.Columns.Item(2).ColumnWidth = 23.68
$newSheet.Columns.Item(3).ColumnWidth = 34.04
$newSheet.Columns.Item(4).ColumnWidth = 40

# --- Populate translation data ---
$newSheet.Range('A1').Value = 'origin'
$newSheet.Range('B1').Value = 'en'
$newSheet.Range('C1').Value = 'ua'
$newSheet.Range('D1').Value = 'ru'

$newSheet.Range('A2').Value = 'PHX_PAGE_NAME'
$newSheet.Range('B2').Value = 'Create arena'
$newSheet.Range('C2').Value = 'Створити арену'
$newSheet.Range('D2').Value = 'Создать арену'

$newSheet.Range('A3').Value = 'PHX_PLACE_NAME'
$newSheet.Range('B3').Value = 'Arena name'
$newSheet.Range('C3').Value = 'Назва арени'
$newSheet.Range('D3').Value = 'Название арены'

$newSheet.Range('A4').Value = 'PHX_AUTOSORT_TAB'
$newSheet.Range('B4').Value = 'Automatic sorting'
$newSheet.Range('C4').Value = 'Автоматичне сортування'
$newSheet.Range('D4').Value = 'Автоматическая сортировка'

$newSheet.Range('A5').Value = 'PHX_MANUAL_TAB'
$newSheet.Range('B5').Value = 'Manual sorting'
$newSheet.Range('C5').Value = 'Ручне сортування'
$newSheet.Range('D5').Value = 'Ручная сортировка'

$newSheet.Range('A6').Value = 'PHX_ACTIVE_GROUPS_LIST'
$newSheet.Range('B6').Value = 'Selected groups'
$newSheet.Range('C6').Value = 'Обрані групи'
$newSheet.Range('D6').Value = 'Выбранные группы'

$newSheet.Range('A7').Value = 'PHX_UNACTIVE_GROUPS_LIST'
$newSheet.Range('B7').Value = 'Groups list'
$newSheet.Range('C7').Value = 'Список груп'
$newSheet.Range('D7').Value = 'Список групп'

$newSheet.Range('A8').Value = 'PHX_DISTANCE'
$newSheet.Range('B8').Value = 'One member pairs distance'
$newSheet.Range('C8').Value = 'Відстань між парами одного учасника'
$newSheet.Range('D8').Value = 'Расстояние между парами одного участника'

$newSheet.Range('A9').Value = 'PHX_AGE_MIN'
$newSheet.Range('B9').Value = 'Minimal age'
$newSheet.Range('C9').Value = 'Мінімальний вік'
$newSheet.Range('D9').Value = 'Минимальный возраст'

$newSheet.Range('A10').Value = 'PHX_AGE_MAX'
$newSheet.Range('B10').Value = 'Maximal age'
$newSheet.Range('C10').Value = 'Максимальний вік'
$newSheet.Range('D10').Value = 'Максимальный возраст'

$newSheet.Range('A11').Value = 'PHX_WEIGHT_MIN'
$newSheet.Range('B11').Value = 'Minimal weight'
$newSheet.Range('C11').Value = 'Мінімальна вага'
$newSheet.Range('D11').Value = 'Минимальный вес'

$newSheet.Range('A12').Value = 'PHX_WEIGHT_MAX'
$newSheet.Range('B12').Value = 'Maximal weight'
$newSheet.Range('C12').Value = 'Максимальна вага'
$newSheet.Range('D12').Value = 'Максимальный вес'

$newSheet.Range('A13').Value = 'PHX_FINAL_MIN'
$newSheet.Range('B13').Value = 'Minimal final part'
$newSheet.Range('C13').Value = 'Мінімальна фінальна частина'
$newSheet.Range('D13').Value = 'Минимальная финальная часть'

$newSheet.Range('A14').Value = 'PHX_FINAL_MAX'
$newSheet.Range('B14').Value = 'Maximal final part'
$newSheet.Range('C14').Value = 'Максимальна фінальна частина'
$newSheet.Range('D14').Value = 'Максимальная финальная часть'

$newSheet.Range('A15').Value = 'PHX_APPLY_BTN'
$newSheet.Range('B15').Value = 'CREATE'
$newSheet.Range('C15').Value = 'СТВОРИТИ'
$newSheet.Range('D15').Value = 'СОЗДАТЬ'

$newSheet.Range('A16').Value = 'PHX_ACTIVE_PAIRS_LIST'
$newSheet.Range('B16').Value = 'Selected pairs'
$newSheet.Range('C16').Value = 'Обрані пари'
$newSheet.Range('D16').Value = 'Выбранные пары'

$newSheet.Range('A17').Value = 'PHX_UNACTIVE_PAIRS_LIST'
$newSheet.Range('B17').Value = 'Unattached pairs'
$newSheet.Range('C17').Value = 'Неприв’язані пари'
$newSheet.Range('D17').Value = 'Непривязанные пары'

# --- Wrap text for the section header cell A2 (style used a wrap-text xf) ---
$newSheet.Range('A2').WrapText = $true

# --- Row heights to roughly match target (cosmetic) ---
$newSheet.StandardHeight = 12.8

# --- Final selection/active cell on the new sheet ---
$newSheet.Activate()
$newSheet.Range('H35').Select()